# Updated symbol list on Sun Feb  5 14:29:15 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (column D) and "Volume(1h)" (column E) quote values
# for the rows whose coin prices/percentages changed in this snapshot.
#
# Both columns store plain text (e.g. "329.34", "-0.58%") rather than
# numbers, so each cell is temporarily switched to Text number format
# before the assignment (and its original style restored right after)
# to stop Excel's automatic number/percentage parsing from mangling the
# literal text (trailing zeros, exact precision, "%" suffix, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-QuoteText($address, $text) {
    $cell = $ws.Range($address)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# row, new Price (column D, $null = unchanged), new Volume(1h) (column E)
$updates = @(
    @{ Row = 2;  D = "329.34";       E = "-0.58%" },
    @{ Row = 3;  D = "43.30";        E = "3.79%"  },
    @{ Row = 4;  D = "5.603";        E = "-2.00%" },
    @{ Row = 5;  D = "0.08196";      E = "-1.72%" },
    @{ Row = 6;  D = "8.795";        E = "-0.02%" },
    @{ Row = 7;  D = "4.402";        E = "-3.11%" },
    @{ Row = 8;  D = "1.902";        E = "-6.44%" },
    @{ Row = 10; D = "0.9433";       E = "1.70%"  },
    @{ Row = 11; D = "0.1201";       E = "-6.38%" },
    @{ Row = 12; D = "0.1921";       E = "-2.20%" },
    @{ Row = 13; D = "0.09815";      E = "4.55%"  },
    @{ Row = 14; D = "0.04356";      E = "11.46%" },
    @{ Row = 15; D = "0.1069";       E = "0.81%"  },
    @{ Row = 16; D = "0.001276";     E = "-2.44%" },
    @{ Row = 17; D = "0.005987";     E = "-2.47%" },
    @{ Row = 18; D = "3.500";        E = "1.72%"  },
    @{ Row = 20; D = "8.739";        E = "6.25%"  },
    @{ Row = 21; D = $null;          E = "-0.15%" },
    @{ Row = 22; D = "0.2522";       E = "4.50%"  },
    @{ Row = 23; D = $null;          E = "-0.66%" },
    @{ Row = 24; D = $null;          E = "-0.74%" },
    @{ Row = 25; D = "0.004301";     E = "-1.80%" },
    @{ Row = 26; D = "0.0001235";    E = "2.81%"  },
    @{ Row = 27; D = "0.0004006";    E = "31.55%" },
    @{ Row = 39; D = "0.02825";      E = "1.07%"  },
    @{ Row = 40; D = "0.05722";      E = "2.75%"  },
    @{ Row = 41; D = "0.007927";     E = "1.69%"  },
    @{ Row = 42; D = "0.009804";     E = "9.59%"  },
    @{ Row = 43; D = $null;          E = "-1.33%" },
    @{ Row = 44; D = "0.002100";     E = "-1.99%" },
    @{ Row = 45; D = "0.01004";      E = "-9.38%" },
    @{ Row = 46; D = "0.00007317";   E = "3.85%"  },
    @{ Row = 47; D = "0.00000000753"; E = "0.31%" },
    @{ Row = 48; D = "0.003496";     E = "-0.34%" },
    @{ Row = 49; D = "0.002280";     E = "-0.03%" },
    @{ Row = 50; D = "0.00002109";   E = "0.31%"  },
    @{ Row = 51; D = $null;          E = "0.31%"  }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-QuoteText "D$($u.Row)" $u.D
    }
    if ($null -ne $u.E) {
        Set-QuoteText "E$($u.Row)" $u.E
    }
}
